# Generate Report for Handback
# Refresh the handoff/handback timestamps recorded for
# "5a779828-e74f-4d8a-b16b-f02ee3ccd84f.md" (row 3 on each sheet) after a
# new localization report was generated.

$wb = $excel.ActiveWorkbook

# Overview sheet: "Latest HO Xliff Generate Date" column (G)
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("G3").Value = "2017-02-09 15:26:27"

# zh-cn sheet: "Correspond Handoff Datetime" (H) / "Correspond Handback DateTime" (L)
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("H3").Value = "2017-02-09 15:26:09"
$wsZhCn.Range("L3").Value = "2017-02-09 15:27:04"

# de-de sheet: "Correspond Handoff Datetime" (H) / "Correspond Handback DateTime" (L)
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("H3").Value = "2017-02-09 15:26:27"
$wsDeDe.Range("L3").Value = "2017-02-09 15:27:30"
